# Change to use sqlalchemy in import database
# - Add a new row of data (row 6) to the "C60 db" sheet describing a new
#   "jq100 devices" batch, using the same date-format style as the rows above.
# - Add a corresponding "Jq101" lookup value to Sheet2 (row 4).
# - Move the active selection/tab from "C60 db" to "Sheet2".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("C60 db")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet "C60 db": fill in row 6 with a new device-batch record ---
$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = 5467
$ws1.Range("C6").Value = "jq100 devices"
$ws1.Range("D6").Value = "Jq100"
$ws1.Range("E6").Value = "Jq150"
# F6 already carries the date-format style from the template; copy the cell
# above it (same "7/2/2012" entry) down so the value/type/style all match.
$ws1.Range("F5").Copy($ws1.Range("F6"))

# --- Sheet "Sheet2": append the new reference value ---
$ws2.Range("A4").Value = "Jq101"

# --- Selections / active sheet, matching the author's final view state ---
$ws1.Range("F9").Select() | Out-Null
$ws2.Range("C7").Select() | Out-Null
$ws2.Activate() | Out-Null
